$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 15

$ws.Cells.Item($row, 1).Value = "mistral:7b-instruct-v0.3-q5_K_M"
$ws.Cells.Item($row, 2).Value = "llama3:70b"
$ws.Cells.Item($row, 3).Value = 1
$ws.Cells.Item($row, 4).Value = 20
$ws.Cells.Item($row, 5).Value = 548.28
$ws.Cells.Item($row, 6).Value = 17.0596
$ws.Cells.Item($row, 7).Value = 100
$ws.Cells.Item($row, 8).Value = "mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_1_20_val.txt"
$ws.Cells.Item($row, 9).Value = 55.76
$ws.Cells.Item($row, 10).Value = 75
$ws.Cells.Item($row, 11).Value = "mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_1_20_test.txt"
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 4
$ws.Cells.Item($row, 14).Value = 20.11
$ws.Cells.Item($row, 15).Value = 100
$ws.Cells.Item($row, 16).Value = "mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_1_20_val_fewshot.txt"
$ws.Cells.Item($row, 17).Value = 66.87
$ws.Cells.Item($row, 18).Value = 75
$ws.Cells.Item($row, 19).Value = "mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_1_20_test_fewshot.txt"
$ws.Cells.Item($row, 20).Value = 302.84
$ws.Cells.Item($row, 21).Value = 2
$ws.Cells.Item($row, 22).Value = 2
$ws.Cells.Item($row, 23).Value = 27
$ws.Cells.Item($row, 24).Value = 100
$ws.Cells.Item($row, 25).Value = "mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_1_20_val_bootstrap.txt"
$ws.Cells.Item($row, 26).Value = 58.64
$ws.Cells.Item($row, 27).Value = 75
$ws.Cells.Item($row, 28).Value = "mistral_7b_instruct_v0.3_q5_K_M_llama3_70b_1_20_test_bootstrap.txt"
